# Reorder / refresh the "Estado de Cuenta" worker table (rows 16-29).
# Previously the rows were grouped by period (all 1707 rows, then all 1708 rows).
# Now they are grouped by worker, each worker's 1708 period followed by their 1707 period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; Doc = "73191450";    Nombre = "OMAR HURTADO MARTINEZ";          Periodo = "1708"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 17; Doc = "73191450";    Nombre = "OMAR HURTADO MARTINEZ";          Periodo = "1707"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 18; Doc = "73007134";    Nombre = "ELKIN JAVIER SALAS CHICO";       Periodo = "1708"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 19; Doc = "73007134";    Nombre = "ELKIN JAVIER SALAS CHICO";       Periodo = "1707"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 20; Doc = "1042453994";  Nombre = "SHADYA MICHELL PALLARES MARIN";  Periodo = "1708"; Valor = 37000;  Salario = 925000 },
    @{ Row = 21; Doc = "1042453994";  Nombre = "SHADYA MICHELL PALLARES MARIN";  Periodo = "1707"; Valor = 37000;  Salario = 925000 },
    @{ Row = 22; Doc = "73152761";    Nombre = "JAVIER ENRIQUE MENDEZ MUÑOZ";    Periodo = "1708"; Valor = 65400;  Salario = 1635000 },
    @{ Row = 23; Doc = "73152761";    Nombre = "JAVIER ENRIQUE MENDEZ MUÑOZ";    Periodo = "1707"; Valor = 65400;  Salario = 1635000 },
    @{ Row = 24; Doc = "10939669";    Nombre = "OSCAR ENRIQUE JULIO BOLAÑOS";    Periodo = "1708"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 25; Doc = "10939669";    Nombre = "OSCAR ENRIQUE JULIO BOLAÑOS";    Periodo = "1707"; Valor = 50140;  Salario = 1253500 },
    @{ Row = 26; Doc = "9090059";     Nombre = "VICTOR MORENO DOMINGUEZ";        Periodo = "1708"; Valor = 58140;  Salario = 1453500 },
    @{ Row = 27; Doc = "9090059";     Nombre = "VICTOR MORENO DOMINGUEZ";        Periodo = "1707"; Valor = 58140;  Salario = 1453500 },
    @{ Row = 28; Doc = "73122616";    Nombre = "ARGEMIRO BARBOZA LUNA";          Periodo = "1708"; Valor = 152600; Salario = 3815000 },
    @{ Row = 29; Doc = "73122616";    Nombre = "ARGEMIRO BARBOZA LUNA";          Periodo = "1707"; Valor = 152600; Salario = 3815000 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("C$i").Value = $r.Doc
    $ws.Range("D$i").Value = $r.Nombre
    $ws.Range("E$i").Value = $r.Periodo
    $ws.Range("F$i").Value = $r.Valor
    $ws.Range("G$i").Value = $r.Salario
}
